$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.960.83'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '1.843.39'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("E5").Value = '  +0.44%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '308.62'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.46%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4769'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.88%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3676'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07199'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.33%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.9289'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.67%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '19.79'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.89%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.07721'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.47%  '
$ws.Range("D13").Value = '1.837.93'
$ws.Range("E13").Value = '  -2.17%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.418'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.36%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.445'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +0.99%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '88.77'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("E17").Value = '  +0.39%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008643'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").Value = '27.025.05'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("E21").Value = '  +1.25%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.068'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("E23").Value = '  -0.28%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.943'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.85%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '152.50'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.01%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '18.18'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.80%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '2.012'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.88%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '114.34'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.26%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '4.962'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.36%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.08859'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.11%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '3.317'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +4.52%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.174'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.61%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.7388'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.45%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '4.491'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.42%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.686'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -5.75%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.107'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.62%  '
$ws.Range("E37").Value = '  +1.06%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.05250'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +1.71%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.963'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.63%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.5246'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +2.26%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '7.002'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("E42").Value = '  -0.23%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '8.275'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.17%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '10.53'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.52%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.4729'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("E46").Value = '  +0.49%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '101.84'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.55%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.601'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.15%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '65.76'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.36%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06068'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.16%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.8862'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.95%  '
